$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant
$xlPasteFormats = -4122
$xlCenter = -4108

# --- 1. Insert two new blank rows before the old row 167 block ---
# (old row167 -> new row169, old row168 -> new row170)
$ws.Rows.Item(167).Insert()
$ws.Rows.Item(167).Insert()

# --- 2. Fix up the merged cell range A167:A168 -> A167:A170 ---
$ws.Range("A169:A170").UnMerge()
$ws.Range("A167:A170").Merge()
# the merge operation resets formatting on the whole range -- restore the
# "centered, no wrap" look shared by the rest of this merged scope column
$ws.Range("A167:A170").HorizontalAlignment = $xlCenter
$ws.Range("A167:A170").WrapText = $false

# --- 3. Copy formats onto the newly inserted rows from the (shifted) old row167 ---
$ws.Range("B169").Copy()
$ws.Range("B167").PasteSpecial($xlPasteFormats)
$ws.Range("B164").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- 4. Row heights for the two newly inserted rows ---
$ws.Rows.Item(167).RowHeight = 27
$ws.Rows.Item(168).RowHeight = 25.5

# --- 5. Row 164: priority bumped to 5, new comment added ---
$ws.Range("B164").Value = 5
$ws.Range("D164").Value = "自动驾驶视觉任务有: 1. object recogonition: 包含具有类别的物体识别和语义分割；2. general object detection: 包括静态物体识别和动态物体识别；3. 距离相关预测， 比如自由度；4. 场景识别，恶劣天气识别和在线校准等；多任务学习不一定能提升准确率，但是能够在减少计算量的前提下获得不错的效果；参考文献很丰富"

# --- 6. Row 167 (new): section header "自动驾驶相关综述" + first reading note ---
$ws.Range("A167").Value = "自动驾驶相关综述"
$ws.Range("B167").Value = 5
$ws.Range("C167").Value = "Autonomous vehicle perception: The technology of today and tomorrow"
$ws.Range("D167").Value = "很好的一篇综述，主要review了无人车上的不同sensor和localizaiton和map的算法，但是没有详细说明preception的任务主要有哪些？"

# --- 7. Row 168 (new): second reading note (no priority, no comment) ---
$ws.Range("C168").Value = "Algorithm and hardware implementation for visual perception system in autonomous vehicle: a survey"

# --- 8. Sheet view bookkeeping: keep the frozen pane / active selection pointed near the bottom ---
$ws.Application.ActiveWindow.ScrollRow = 158
$ws.Range("C164").Select()
